$d = $word.ActiveDocument

# 1. Rewrite the opening sentence (before the first "SelfHelp").
$d.Content.Find.Execute(
    "Developing an app requires an idea, programming skills, and time. For non-developers, app development can be daunting. However, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "App development requires an idea, programming skills, and time, which can be daunting for non-developers. ",
    2) | Out-Null

# 2. Between the first and second "SelfHelp" mentions.
$d.Content.Find.Execute(
    ", an open-source CMS system eliminates the need for programming skills, offering components designed to meet researchers' needs. With ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " is an open-source CMS system that eliminates the need for programming skills, offering components designed to meet researchers' needs. With ",
    2) | Out-Null

# 3. Between the second and third "SelfHelp" mentions: reorder sentences.
$d.Content.Find.Execute(
    ", users can design pages, menus, and forms without requiring programming skills. Its interactive data collection enables dynamic event triggers and its application functionality to be customized based on the incoming data, making it easy to design intervention studies with multiple sessions and pre- and post-tests. Condition components can be used to fine-tune what is displayed, when it is displayed, and how long it remains accessible. Additionally, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", users can design pages, menus, and forms without requiring programming skills. Condition components can be used to fine-tune what is displayed, when it is displayed, and how long it remains accessible. By utilizing the incoming data and the scheduling system it is easy to design intervention studies with multiple sessions and pre- and post-tests. Additionally, ",
    2) | Out-Null

Write-Output ("Para2 = [" + $d.Paragraphs(2).Range.Text + "]")
